# Remove the empty paragraph and the "Vi bedömer..." paragraph that
# followed the "Uppgången är också..." paragraph in the CEO comment
# section. After the edit, "Uppgången är också..." is immediately
# followed by the blank paragraph that used to precede
# "Rörelseresultatet steg...".

$d = $word.ActiveDocument

$marker = "Vi bedömer att bortfallet i omsättning"

$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*$marker*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -lt 0) {
    throw "Could not locate the 'Vi bedömer...' paragraph."
}

# The paragraph immediately before it is the blank spacer paragraph;
# remove both of them as a single range so the surrounding paragraphs
# collapse back together cleanly.
$target = $d.Paragraphs.Item($targetIndex)
$prev = $d.Paragraphs.Item($targetIndex - 1)

$startPos = $prev.Range.Start
$endPos = $target.Range.End

$r = $d.Range($startPos, $endPos)
$r.Delete()
